# #5: property boat&car done
# Adds the "capacity" header column plus the common
# property_category/category/date/legislator_name/legislator_id/source_file/index
# columns (H:N) to the 汽車 (car) sheet, matching the layout already used on
# the other property sheets (土地, 股票, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: turn the old "sample data" header row into a real header row ---
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# match the bold/border header formatting already used on B1:G1
$ws.Cells.Item(1, 2).Copy()
$ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 14)).PasteSpecial(-4122)

# --- Rows 2-4: fill in the new property/legislator metadata columns ---
$carRows = @(2, 3, 4)
foreach ($r in $carRows) {
    $ws.Cells.Item($r, 8).Value = "land"
    $ws.Cells.Item($r, 9).Value = "normal"
    # leading apostrophe forces text (otherwise Excel reparses this as a date serial)
    $ws.Cells.Item($r, 10).Value = "'2012-03-06"
    $ws.Cells.Item($r, 11).Value = "廖國棟"
    $ws.Cells.Item($r, 12).Value = 962
    $ws.Cells.Item($r, 13).Value = "tmpec731"
}
$ws.Cells.Item(2, 14).Value = 33
$ws.Cells.Item(3, 14).Value = 34
$ws.Cells.Item(4, 14).Value = 37

# match the plain body formatting already used on the existing data cells
$ws.Cells.Item(2, 2).Copy()
$ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item(4, 14)).PasteSpecial(-4122)
